# Auto-update draw results: append the newest Pick 3 draw as a new row
# at the bottom of the "Results" table (row 58), mirroring the existing
# rows above it (columns: Date, Game, Phase, Result, InsertedAt).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 58

# Columns A (Date) and C (Phase) look like a date / a plain number, so
# Excel would otherwise silently convert them on entry. Force text
# formatting on just those two cells first so the stored values stay
# literal strings, exactly like every other row in the table.
$ws.Range("A" + $row).NumberFormat = "@"
$ws.Range("C" + $row).NumberFormat = "@"

$ws.Range("A" + $row).Value = "2025-11-13"
$ws.Range("B" + $row).Value = "Pick 3"
$ws.Range("C" + $row).Value = "251113"
$ws.Range("D" + $row).Value = "6-8-5"
$ws.Range("E" + $row).Value = "2025-11-13T21:41:09.310+04:00"
